# Update scripts wuth new tpm
# Recomputed ligand/receptor/edge expression-specificity metrics (columns
# G:T) for each sending/target cluster pair after re-running the NATMI TPM
# pipeline. Only the cells whose derived values actually moved are touched;
# identifier/category columns (A:F) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G=3.685507; H=11.056521; I=0.3585631737883472; J=0.3585631737883472; K=3; L=1; M=1.572065666666667; N=4.716197; O=0.1759712293834306; P=0.1759712293834305; Q=5.793859018959667; R=52.144731170637; S=0.06309680250316013; T=0.06309680250316012 }
    3  = @{ G=3.685507; H=11.056521; I=0.3585631737883472; J=0.3585631737883472; O=0.4743638053196239; P=0.4743638053196239; Q=15.618445250107; R=140.566007250963; S=0.170089391565722; T=0.170089391565722 }
    4  = @{ G=3.685507; H=11.056521; I=0.3585631737883472; J=0.3585631737883472; M=3.123785; N=9.371354999999999; O=0.3496649652969456; P=0.3496649652969455; Q=11.512731483995; R=103.614583355955; S=0.1253769797194651; T=0.1253769797194651 }
    5  = @{ I=0.009647184430711629; J=0.009647184430711629; K=3; L=1; M=1.572065666666667; N=4.716197; O=0.1759712293834306; P=0.1759712293834305; Q=0.155884459441; R=1.402960134969; S=0.001697626904361016; T=0.001697626904361016 }
    6  = @{ I=0.009647184430711629; J=0.009647184430711629; O=0.4743638053196239; P=0.4743638053196239; S=0.004576275117172598; T=0.004576275117172598 }
    7  = @{ I=0.009647184430711629; J=0.009647184430711629; M=3.123785; N=9.371354999999999; O=0.3496649652969456; P=0.3496649652969455; Q=0.309751396815; R=2.787762571335; S=0.003373282409178015; T=0.003373282409178015 }
    8  = @{ G=6.493877; H=19.481631; I=0.6317896417809412; J=0.6317896417809411; K=3; L=1; M=1.572065666666667; N=4.716197; O=0.1759712293834306; P=0.1759712293834305; Q=10.20880107525634; R=91.879209677307; S=0.1111767999759094; T=0.1111767999759094 }
    9  = @{ G=6.493877; H=19.481631; I=0.6317896417809412; J=0.6317896417809411; O=0.4743638053196239; P=0.4743638053196239; Q=27.519758444477; R=247.677826000293; S=0.2996981386367293; T=0.2996981386367293 }
    10 = @{ G=6.493877; H=19.481631; I=0.6317896417809412; J=0.6317896417809411; M=3.123785; N=9.371354999999999; O=0.3496649652969456; P=0.3496649652969455; Q=20.285475564445; R=182.569280080005; S=0.2209147031683025; T=0.2209147031683024 }
}

foreach ($row in $updates.Keys) {
    $rowVals = $updates[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
